# Apply the "Updated cryptos list" data refresh to Sheet1.
# Rows 2-51 hold the scraped crypto table (Coin, Link, Price, Volume(1h));
# this updates the Price/Volume figures, and rows 33-34 additionally swap
# which coin (Bittensor vs FirstDigitalUSD) occupies each rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.645.07"
$ws.Range("E2").Value = "  -1.25%  "
# Row 3
$ws.Range("D3").Value = "2.450.19"
$ws.Range("E3").Value = "  -1.60%  "
# Row 4
$ws.Range("E4").Value = "  -0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.36%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "
# Row 7
$ws.Range("E7").Value = "  -0.06%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.83%  "
# Row 9
$ws.Range("D9").Value = "2.449.57"
$ws.Range("E9").Value = "  -1.55%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.71%  "
# Row 11
$ws.Range("E11").Value = "  -0.90%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.96%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "
# Row 14
$ws.Range("D14").Value = "2.890.50"
$ws.Range("E14").Value = "  -1.89%  "
# Row 15
$ws.Range("D15").Value = "68.418.26"
$ws.Range("E15").Value = "  -1.41%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.16%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.52%  "
# Row 18
$ws.Range("D18").Value = "2.407.25"
$ws.Range("E18").Value = "  -3.73%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.97%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.07%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.78%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.37%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "
# Row 24
$ws.Range("E24").Value = "  +0.01%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.32%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.17%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.93%  "
# Row 28
$ws.Range("D28").Value = "2.562.44"
$ws.Range("E28").Value = "  -2.09%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.83%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0812"
$ws.Range("E31").Value = "  -6.61%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.79%  "
# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "434.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
# Row 35
$ws.Range("E35").Value = "  -5.40%  "
# Row 36
$ws.Range("E36").Value = "  -5.82%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
# Row 39
$ws.Range("E39").Value = "  +0.12%  "
# Row 40
$ws.Range("E40").Value = "  -3.42%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.56%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.03%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.57%  "
# Row 46
$ws.Range("E46").Value = "  +1.73%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.58%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.85%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0712"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.480"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.74%  "
